$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "barnehage": rename kindergartens to real Drammen Kommune names
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("barnehage")

$wsB.Range("C2").Value = "Amicus barnehage"
$wsB.Range("C3").Value = "Aronsløkka barnehage"
$wsB.Range("C4").Value = "Bacheparken barnehage"
$wsB.Range("C5").Value = "Blåbærtoppen barnehage"
$wsB.Range("C6").Value = "Dalegårdsveien barnehage"
$wsB.Range("C7").Value = "Ekornheia barnehage"
$wsB.Range("D7").Value = 20
$wsB.Range("C8").Value = "Eplehagan barnehage"

# ---------------------------------------------------------------------------
# Sheet "soknad": update existing rows + add new rows with new column Q
# ---------------------------------------------------------------------------
$wsS = $wb.Worksheets.Item("soknad")

# New header for column Q
$wsS.Range("Q1").Value = "fortrinnsrett"

# Row 2 (existing, updated)
$wsS.Range("N2").Value = "Haakon Halvorsen"
$wsS.Range("O2").Value = "Amicus barnehage"
$wsS.Range("P2").Value = "TILBUD"
$wsS.Range("Q2").Value = "Nei"

# Row 3 (existing, updated)
$wsS.Range("N3").Value = "Haakon Halvorsen"
$wsS.Range("O3").Value = "Amicus barnehage"
$wsS.Range("P3").Value = "TILBUD"
$wsS.Range("Q3").Value = "Ja"

# Row 4 (existing, updated)
$wsS.Range("N4").Value = "Haakon Halvorsen"
$wsS.Range("O4").Value = "Blåbærtoppen barnehage"
$wsS.Range("P4").Value = "AVSLAG: Ingen ledige plasser."
$wsS.Range("Q4").Value = "Nei"

# Row 5 (new)
$wsS.Range("N5").Value = "Haakon Halvorsen"
$wsS.Range("O5").Value = "Blåbærtoppen barnehage"
$wsS.Range("P5").Value = "AVSLAG: Ingen ledige plasser."
$wsS.Range("Q5").Value = "Nei"

# Row 6 (new)
$wsS.Range("N6").Value = "Haakon Halvorsen"
$wsS.Range("O6").Value = "Blåbærtoppen barnehage"
$wsS.Range("P6").Value = "AVSLAG: Ingen ledige plasser."
$wsS.Range("Q6").Value = "Nei"

# Row 7 (new)
$wsS.Range("N7").Value = "Haakon Halvorsen"
$wsS.Range("O7").Value = "Amicus barnehage"
$wsS.Range("P7").Value = "TILBUD"
$wsS.Range("Q7").Value = "Nei"

# Row 8 (new)
$wsS.Range("N8").Value = "Haakon Halvorsen"
$wsS.Range("O8").Value = "Amicus barnehage"
$wsS.Range("P8").Value = "TILBUD"
$wsS.Range("Q8").Value = "Nei"

# Row 9 (new)
$wsS.Range("N9").Value = "Haakon Halvorsen"
$wsS.Range("O9").Value = "Eplehagan barnehage"
$wsS.Range("P9").Value = "TILBUD"
$wsS.Range("Q9").Value = "Nei"

# Row 10 (new)
$wsS.Range("N10").Value = "Haakon Halvorsen"
$wsS.Range("O10").Value = "Eplehagan barnehage"
$wsS.Range("P10").Value = "TILBUD"
$wsS.Range("Q10").Value = "Nei"

# Row 11 (new)
$wsS.Range("N11").Value = "Ola Nordmann"
$wsS.Range("O11").Value = "Blåbærtoppen barnehage"
$wsS.Range("P11").Value = "AVSLAG: Barnet er under ett år."
$wsS.Range("Q11").Value = "Nei"

# Row 12 (new)
$wsS.Range("N12").Value = "Ola Nordmann"
$wsS.Range("O12").Value = "Blåbærtoppen barnehage"
$wsS.Range("P12").Value = "AVSLAG: Ingen ledige plasser."
$wsS.Range("Q12").Value = "Nei"

# Row 13 (new)
$wsS.Range("N13").Value = "Ola Nordmann"
$wsS.Range("O13").Value = "Amicus barnehage"
$wsS.Range("P13").Value = "TILBUD"
$wsS.Range("Q13").Value = "Nei"

# Row 14 (new)
$wsS.Range("N14").Value = "Ola Nordmann"
$wsS.Range("O14").Value = "Amicus barnehage"
$wsS.Range("P14").Value = "TILBUD"
$wsS.Range("Q14").Value = "Nei"

# Row 15 (new)
$wsS.Range("N15").Value = "Ola Nordmann"
$wsS.Range("O15").Value = "Amicus barnehage"
$wsS.Range("P15").Value = "TILBUD"
$wsS.Range("Q15").Value = "Nei"

# Row 16 (new)
$wsS.Range("N16").Value = "Ola Nordmann"
$wsS.Range("O16").Value = "Amicus barnehage"
$wsS.Range("P16").Value = "TILBUD"
$wsS.Range("Q16").Value = "Nei"

# Row 17 (new)
$wsS.Range("N17").Value = "Ola Nordmann"
$wsS.Range("O17").Value = "Amicus barnehage"
$wsS.Range("P17").Value = "TILBUD"
$wsS.Range("Q17").Value = "Nei"

# Row 18 (new)
$wsS.Range("N18").Value = "Ola Nordmann"
$wsS.Range("O18").Value = "Amicus barnehage"
$wsS.Range("P18").Value = "TILBUD"
$wsS.Range("Q18").Value = "Nei"

# Row 19 (new) - Q19 left empty (no preference value)
$wsS.Range("N19").Value = "Ola Nordmann"
$wsS.Range("O19").Value = "Dalegårdsveien barnehage"
$wsS.Range("P19").Value = "TILBUD"

# Row 20 (new) - Q20 left empty (no preference value)
$wsS.Range("N20").Value = "Ola Nordmann"
$wsS.Range("O20").Value = "Dalegårdsveien barnehage"
$wsS.Range("P20").Value = "AVSLAG: Barnet er under ett år."
